# Regenerate save_data to use K (column G) instead of Strike#.
# This updates the previously-pulled "K" values in column G for each
# data row with freshly recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 0
    21 = 0
    22 = 1
    23 = 3
    24 = 1
    25 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
